# Apply updated "Price" (D) and "Volume(1h)" (E) values scraped for this run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.125.64'
$ws.Range('E2').Value = '  +0.92%  '
$ws.Range('D3').Value = '1.639.05'
$ws.Range('E3').Value = '  -0.07%  '
$ws.Range('E4').Value = '  +0.06%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = "@"
$cell.Value = '216.70'
$cell.Style = "Normal"
$ws.Range('E5').Value = '  -0.09%  '
$ws.Range('E6').Value = '  +1.81%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('E8').Value = '  -0.26%  '
$ws.Range('E9').Value = '  -0.02%  '
$cell = $ws.Range('D10')
$cell.NumberFormat = "@"
$cell.Value = '19.98'
$cell.Style = "Normal"
$ws.Range('E10').Value = '  +0.41%  '
$cell = $ws.Range('D11')
$cell.NumberFormat = "@"
$cell.Value = '0.0848'
$cell.Style = "Normal"
$ws.Range('E11').Value = '  +0.11%  '
$ws.Range('D12').Value = '1.868.27'
$ws.Range('E12').Value = '  -0.08%  '
$ws.Range('D13').Value = '1.636.80'
$ws.Range('E13').Value = '  +0.17%  '
$ws.Range('E14').Value = '  +0.10%  '
$ws.Range('E15').Value = '  +1.94%  '
$cell = $ws.Range('D16')
$cell.NumberFormat = "@"
$cell.Value = '66.73'
$cell.Style = "Normal"
$ws.Range('E16').Value = '  -1.10%  '
$ws.Range('D17').Value = '27.131.00'
$ws.Range('E17').Value = '  +0.98%  '
$ws.Range('E18').Value = '  +0.95%  '
$cell = $ws.Range('D19')
$cell.NumberFormat = "@"
$cell.Value = '216.97'
$cell.Style = "Normal"
$ws.Range('E19').Value = '  -1.18%  '
$ws.Range('E20').Value = '  +0.03%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = "@"
$cell.Value = '6.97'
$cell.Style = "Normal"
$ws.Range('E21').Value = '  +1.82%  '
$ws.Range('E22').Value = '  +3.61%  '
$ws.Range('E23').Value = '  +0.33%  '
$ws.Range('E24').Value = '  -0.31%  '
$cell = $ws.Range('D25')
$cell.NumberFormat = "@"
$cell.Value = '146.81'
$cell.Style = "Normal"
$ws.Range('E25').Value = '  -0.38%  '
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('E27').Value = '  +1.01%  '
$ws.Range('E28').Value = '  +0.06%  '
$cell = $ws.Range('D29')
$cell.NumberFormat = "@"
$cell.Value = '15.68'
$cell.Style = "Normal"
$ws.Range('E29').Value = '  -0.85%  '
$ws.Range('E30').Value = '  +1.05%  '
$ws.Range('E31').Value = '  -0.23%  '
$ws.Range('E32').Value = '  +1.28%  '
$cell = $ws.Range('D33')
$cell.NumberFormat = "@"
$cell.Value = '3.01'
$cell.Style = "Normal"
$ws.Range('E33').Value = '  +0.48%  '
$ws.Range('D34').Value = '1.307.96'
$ws.Range('E34').Value = '  +2.95%  '
$ws.Range('E35').Value = '  -0.16%  '
$ws.Range('E36').Value = '  +1.28%  '
$ws.Range('E37').Value = '  -1.35%  '
$cell = $ws.Range('D38')
$cell.NumberFormat = "@"
$cell.Value = '0.858'
$cell.Style = "Normal"
$ws.Range('E38').Value = '  +2.80%  '
$cell = $ws.Range('D39')
$cell.NumberFormat = "@"
$cell.Value = '0.542'
$cell.Style = "Normal"
$ws.Range('E39').Value = '  +1.43%  '
$ws.Range('E40').Value = '  +0.02%  '
$cell = $ws.Range('D41')
$cell.NumberFormat = "@"
$cell.Value = '0.811'
$cell.Style = "Normal"
$ws.Range('E41').Value = '  +0.27%  '
$ws.Range('E42').Value = '  +5.50%  '
$ws.Range('D44').Value = '1.778.57'
$ws.Range('E44').Value = '  -0.12%  '
$cell = $ws.Range('D45')
$cell.NumberFormat = "@"
$cell.Value = '61.83'
$cell.Style = "Normal"
$ws.Range('E45').Value = '  -0.30%  '
$cell = $ws.Range('D46')
$cell.NumberFormat = "@"
$cell.Value = '91.41'
$cell.Style = "Normal"
$ws.Range('E46').Value = '  -0.36%  '
$ws.Range('E47').Value = '  +0.78%  '
$ws.Range('D48').Value = '0.0₆0108'
$ws.Range('E48').Value = '  +2.40%  '
$ws.Range('E49').Value = '  -0.08%  '
$cell = $ws.Range('D50')
$cell.NumberFormat = "@"
$cell.Value = '7.67'
$cell.Style = "Normal"
$ws.Range('E50').Value = '  +0.35%  '
$cell = $ws.Range('D51')
$cell.NumberFormat = "@"
$cell.Value = '0.0960'
$cell.Style = "Normal"
$ws.Range('E51').Value = '  -0.16%  '
